$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.525.96"
$ws.Range("E2").Value = "  -6.81%  "

# Row 3
$ws.Range("D3").Value = "2.920.43"
$ws.Range("E3").Value = "  -9.29%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -10.02%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.22"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -13.82%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.30%  "

# Row 8
$ws.Range("D8").Value = "2.901.43"
$ws.Range("E8").Value = "  -9.67%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.461"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -15.68%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.143"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -18.42%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -11.40%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.430"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -13.87%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "31.92"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -18.69%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000201"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -18.42%  "

# Row 15
$ws.Range("D15").Value = "3.381.60"
$ws.Range("E15").Value = "  -9.57%  "

# Row 16
$ws.Range("D16").Value = "62.369.22"
$ws.Range("E16").Value = "  -7.02%  "

# Row 17
$ws.Range("E17").Value = "  -5.68%  "

# Row 18
$ws.Range("D18").Value = "2.919.28"
$ws.Range("E18").Value = "  -9.32%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "466.03"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -13.25%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -14.71%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.78"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -14.66%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.630"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -17.63%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.51"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -18.34%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.75"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -12.88%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.84"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -15.22%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.18%  "

# Row 27
$ws.Range("E27").Value = "  -20.10%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.89"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -14.56%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.97"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -15.12%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "24.34"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -17.07%  "

# Row 31
$ws.Range("B31").Value = "Mantle"
$ws.Range("C31").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.04"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -9.41%  "

# Row 32
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.34"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -12.54%  "

# Row 33
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.42%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "470.36"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -14.60%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.27"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -6.17%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.42"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -17.35%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.74"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -17.61%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0386"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -10.25%  "

# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.114"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -9.93%  "

# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0743"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -14.58%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.78"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -17.61%  "

# Row 42
$ws.Range("D42").Value = "2.681.88"
$ws.Range("E42").Value = "  -8.54%  "

# Row 43
$ws.Range("E43").Value = "  -0.21%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.24"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -16.70%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.221"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -16.89%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "112.37"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -8.67%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.101"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -11.53%  "

# Row 48
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.80"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -16.17%  "

# Row 49
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").Value = "0.0₃0469"
$ws.Range("E49").Value = "  -20.85%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.73"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -18.81%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.17"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.20%  "
